# Fix: added missing argument metrics_dic_d2d in bs2d_ul_outputs() call
# Update boolean flags on row 2 of Sheet1: uplink, save_scenario_xlsx,
# save_metrics_xlsx, and show_video should be TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = $true   # uplink
$ws.Range("K2").Value = $true   # save_scenario_xlsx
$ws.Range("L2").Value = $true   # save_metrics_xlsx
$ws.Range("M2").Value = $true   # show_video
